# "version final sin errores"
#
# Changes applied to the "Metadata" sheet of the ValueSet workbook:
#   1. Bump the Version value from 0.4.0 to 0.7.0.
#   2. Remove the Jurisdiction / Chile row entirely (rows below shift up).
#
# The "Include from Precauciones Hos" sheet's visible content is unchanged;
# it only shifts because shared strings are renumbered after the deletion,
# which Excel (and this engine) handles automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Update the Version property value (row 3, column B).
$ws.Cells.Item(3, 2).Value = "0.7.0"

# 2. Delete the entire "Jurisdiction" / "Chile" row (row 11).
$ws.Rows.Item(11).Delete()
